$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows above row 162 (shifts old rows 162:238 down to 168:244),
# duplicating the formatting of row 162 (carries the date-format style on column D).
$ws.Range("A162:R167").EntireRow.Insert()

# Common (unchanging) values copied from the template row block.
$mercadoId = 11
$mercado = "Vega Monumental Concepción"
$region = "Bíobío"
$codreg = 8
$categoriaId = 100112027
$categoria = "Melón"
$unidad = "`$/unidad"
$kgOUnidades = 1
$clasificacion = "Hortaliza"

# New weekly price block (fecha serial 44609) for "Región de O'Higgins".
$newRows = @(
    @{ Row = 162; Variedad = "Calameño"; Calidad = "Extra";   Fecha = 44609; Volumen = 1000; PrecioMin = 800; PrecioMax = 800; PrecioProm = 800; Origen = "Región de O'Higgins"; PrecioKg = 800 },
    @{ Row = 163; Variedad = "Calameño"; Calidad = "Primera"; Fecha = 44609; Volumen = 1000; PrecioMin = 700; PrecioMax = 700; PrecioProm = 700; Origen = "Región de O'Higgins"; PrecioKg = 700 },
    @{ Row = 164; Variedad = "Calameño"; Calidad = "Segunda"; Fecha = 44609; Volumen = 1000; PrecioMin = 600; PrecioMax = 600; PrecioProm = 600; Origen = "Región de O'Higgins"; PrecioKg = 600 },
    @{ Row = 165; Variedad = "Tuna";     Calidad = "Extra";   Fecha = 44609; Volumen = 1000; PrecioMin = 800; PrecioMax = 800; PrecioProm = 800; Origen = "Región de O'Higgins"; PrecioKg = 800 },
    @{ Row = 166; Variedad = "Tuna";     Calidad = "Primera"; Fecha = 44609; Volumen = 1000; PrecioMin = 700; PrecioMax = 700; PrecioProm = 700; Origen = "Región de O'Higgins"; PrecioKg = 700 },
    @{ Row = 167; Variedad = "Tuna";     Calidad = "Segunda"; Fecha = 44609; Volumen = 1000; PrecioMin = 600; PrecioMax = 600; PrecioProm = 600; Origen = "Región de O'Higgins"; PrecioKg = 600 }
)

foreach ($rd in $newRows) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $rd.Fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $categoriaId
    $ws.Cells.Item($r, 7).Value = $categoria
    $ws.Cells.Item($r, 8).Value = $rd.Variedad
    $ws.Cells.Item($r, 9).Value = $rd.Calidad
    $ws.Cells.Item($r, 10).Value = $rd.Volumen
    $ws.Cells.Item($r, 11).Value = $rd.PrecioMin
    $ws.Cells.Item($r, 12).Value = $rd.PrecioMax
    $ws.Cells.Item($r, 13).Value = $rd.PrecioProm
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $rd.Origen
    $ws.Cells.Item($r, 16).Value = $rd.PrecioKg
    $ws.Cells.Item($r, 17).Value = $kgOUnidades
    $ws.Cells.Item($r, 18).Value = $clasificacion
}
